# Generate Report for Handoff
# Refreshes the handoff/handback status + timestamps for the report and
# tightens the now-shorter "Status" columns that used to hold the long
# "Handed back: in sync with en-US" text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-11-15 17:34:17"

# Status column text got shorter, so the column shrinks back down.
$wsOverview.Range("E1").ColumnWidth = 16.3
$wsOverview.Range("F1").ColumnWidth = 16.3

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-11-15 17:34:01"
$wsZhCn.Range("C1").ColumnWidth = 16.3

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-11-15 17:34:17"
$wsDeDe.Range("C1").ColumnWidth = 16.3
